$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 146.5
$ws.Range("I12").Value = 146.5
$ws.Range("K12").Value = 146.5
$ws.Range("M12").Value = 23.5
$ws.Range("H19").Value = 50003572
$ws.Range("I19").Value = 5274.1665
$ws.Range("J19").Value = 125001020
$ws.Range("K19").Value = 5274.1665
$ws.Range("L19").Value = 125001020
$ws.Range("M19").Value = -5099.1665
$ws.Range("N19").Value = -125001370
$ws.Range("H62").Value = 6375.125
$ws.Range("I62").Value = 1999.3334
$ws.Range("K62").Value = 1999.3334
$ws.Range("M62").Value = -1375.3334
$ws.Range("H65").Value = 6375.125
$ws.Range("I65").Value = 1999.3334
$ws.Range("K65").Value = 9996.666999999999
$ws.Range("M65").Value = -6876.666999999999
$ws.Range("H74").Value = 6353.9165
$ws.Range("I74").Value = 5370.6665
$ws.Range("K74").Value = 5370.6665
$ws.Range("M74").Value = -4434.6665
$ws.Range("H77").Value = 6353.9165
$ws.Range("I77").Value = 5370.6665
$ws.Range("K77").Value = 26853.3325
$ws.Range("M77").Value = -22173.3325
$ws.Range("H86").Value = 4118.9355
$ws.Range("I86").Value = 3720.1
$ws.Range("J86").Value = 4308.857
$ws.Range("K86").Value = 3720.1
$ws.Range("L86").Value = 4308.857
$ws.Range("M86").Value = -2597.1
$ws.Range("N86").Value = -6554.857
$ws.Range("H89").Value = 4118.9355
$ws.Range("I89").Value = 3720.1
$ws.Range("J89").Value = 4308.857
$ws.Range("K89").Value = 18600.5
$ws.Range("L89").Value = 21544.285
$ws.Range("M89").Value = -12984.5
$ws.Range("N89").Value = -32776.285
$ws.Range("H98").Value = 510718.25
$ws.Range("I98").Value = 1266.8334
$ws.Range("J98").Value = 1020169.7
$ws.Range("K98").Value = 1266.8334
$ws.Range("L98").Value = 1020169.7
$ws.Range("M98").Value = 231.1666
$ws.Range("N98").Value = -1023165.7
$ws.Range("H106").Value = 9155.125
$ws.Range("I106").Value = 5634.077
$ws.Range("K106").Value = 5634.077
$ws.Range("M106").Value = -5003.077
$ws.Range("H122").Value = 510718.25
$ws.Range("I122").Value = 1266.8334
$ws.Range("J122").Value = 1020169.7
$ws.Range("K122").Value = 3800.5002
$ws.Range("L122").Value = 3060509.1
$ws.Range("M122").Value = -1350.5002
$ws.Range("N122").Value = -3065409.1
$ws.Range("H129").Value = 2217.9375
$ws.Range("I129").Value = 1663.25
$ws.Range("J129").Value = 2772.625
$ws.Range("K129").Value = 4989.75
$ws.Range("L129").Value = 8317.875
$ws.Range("M129").Value = 10.25
$ws.Range("N129").Value = -18317.875
$ws.Range("H131").Value = 2496.3845
$ws.Range("I131").Value = 729.5789
$ws.Range("J131").Value = 7292
$ws.Range("K131").Value = 2188.7367
$ws.Range("L131").Value = 21876
$ws.Range("M131").Value = 2851.2633
$ws.Range("N131").Value = -31956
$ws.Range("H132").Value = 2547.75
$ws.Range("I132").Value = 2274.3635
$ws.Range("J132").Value = 5555
$ws.Range("K132").Value = 6823.0905
$ws.Range("L132").Value = 16665
$ws.Range("M132").Value = -4293.0905
$ws.Range("N132").Value = -21725
$ws.Range("H138").Value = 3240.5535
$ws.Range("J138").Value = 4378.5
$ws.Range("L138").Value = 13135.5
$ws.Range("N138").Value = -23415.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3769.2622
$ws.Range("I32").Value = 3549.8909
$ws.Range("K32").Value = 3549.8909
$ws.Range("M32").Value = -3262.8909
$ws.Range("H122").Value = 35716916
$ws.Range("I122").Value = 1138.5294
$ws.Range("K122").Value = 3415.5882
$ws.Range("M122").Value = -965.5881999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2567.5
$ws.Range("I86").Value = 2113.8667
$ws.Range("K86").Value = 2113.8667
$ws.Range("M86").Value = -990.8667
$ws.Range("H89").Value = 2567.5
$ws.Range("I89").Value = 2113.8667
$ws.Range("K89").Value = 10569.3335
$ws.Range("M89").Value = -4953.333500000001
$ws.Range("H107").Value = 4382.25
$ws.Range("I107").Value = 6442.6665
$ws.Range("K107").Value = 6442.6665
$ws.Range("M107").Value = -4522.6665
$ws.Range("H134").Value = 2270.6155
$ws.Range("I134").Value = 1683.909
$ws.Range("K134").Value = 5051.727000000001
$ws.Range("M134").Value = -2516.727000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12564.8
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("H86").Value = 7000
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877
$ws.Range("H89").Value = 7000
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384
$ws.Range("H105").Value = 4247.5
$ws.Range("I105").Value = 1994.8334
$ws.Range("K105").Value = 1994.8334
$ws.Range("M105").Value = -247.8334
$ws.Range("H134").Value = 10671
$ws.Range("J134").Value = 10671
$ws.Range("L134").Value = 32013
$ws.Range("N134").Value = -37083
$ws.Range("H136").Value = 12564.8
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1314.3334
$ws.Range("I92").Value = 498.6
$ws.Range("K92").Value = 1495.8
$ws.Range("M92").Value = -247.8000000000002
$ws.Range("H99").Value = 10857.75
$ws.Range("I99").Value = 9209.799999999999
$ws.Range("K99").Value = 27629.4
$ws.Range("M99").Value = -25383.4
$ws.Range("H100").Value = 10332.429
$ws.Range("I100").Value = 1999.5
$ws.Range("K100").Value = 5998.5
$ws.Range("M100").Value = -5187.5
$ws.Range("H129").Value = 6947329.5
$ws.Range("J129").Value = 20840592
$ws.Range("L129").Value = 62521776
$ws.Range("N129").Value = -62531776
$ws.Range("H140").Value = 2505.1765
$ws.Range("I140").Value = 2066.2666
$ws.Range("K140").Value = 6198.7998
$ws.Range("M140").Value = -1018.7998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2670207.8
$ws.Range("I102").Value = 3511852.2
$ws.Range("K102").Value = 3511852.2
$ws.Range("M102").Value = -3510230.2
$ws.Range("H107").Value = 1308.4286
$ws.Range("I107").Value = 696.63635
$ws.Range("J107").Value = 3551.6667
$ws.Range("K107").Value = 696.63635
$ws.Range("L107").Value = 3551.6667
$ws.Range("M107").Value = 1223.36365
$ws.Range("N107").Value = -7391.6667
$ws.Range("H132").Value = 3297.0417
$ws.Range("I132").Value = 2524.75
$ws.Range("J132").Value = 7158.5
$ws.Range("K132").Value = 7574.25
$ws.Range("L132").Value = 21475.5
$ws.Range("M132").Value = -5044.25
$ws.Range("N132").Value = -26535.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1000000000
$ws.Range("I13").Value = 1000000000
$ws.Range("K13").Value = 1000000000
$ws.Range("M13").Value = -999999860
$ws.Range("H61").Value = 3647.853
$ws.Range("I61").Value = 2608.96
$ws.Range("K61").Value = 2608.96
$ws.Range("M61").Value = -2406.96
$ws.Range("H82").Value = 13256.134
$ws.Range("I82").Value = 20333.334
$ws.Range("K82").Value = 20333.334
$ws.Range("M82").Value = -19972.334
$ws.Range("H85").Value = 13256.134
$ws.Range("I85").Value = 20333.334
$ws.Range("K85").Value = 20333.334
$ws.Range("M85").Value = -19085.334
$ws.Range("H113").Value = 3647.853
$ws.Range("I113").Value = 2608.96
$ws.Range("K113").Value = 2608.96
$ws.Range("M113").Value = -438.96
$ws.Range("H135").Value = 54253
$ws.Range("J135").Value = 54253
$ws.Range("L135").Value = 54253
$ws.Range("N135").Value = -64393
$ws.Range("H136").Value = 5461.3667
$ws.Range("I136").Value = 4422.846
$ws.Range("K136").Value = 13268.538
$ws.Range("M136").Value = -10718.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1414.7858
$ws.Range("I113").Value = 2359.6
$ws.Range("J113").Value = 889.8889
$ws.Range("K113").Value = 7078.799999999999
$ws.Range("L113").Value = 2669.6667
$ws.Range("M113").Value = -4908.799999999999
$ws.Range("N113").Value = -7009.6667
$ws.Range("H122").Value = 3756.7
$ws.Range("J122").Value = 5084
$ws.Range("L122").Value = 15252
$ws.Range("N122").Value = -20152
$ws.Range("H132").Value = 2214.59
$ws.Range("I132").Value = 1284.6938
$ws.Range("K132").Value = 3854.0814
$ws.Range("M132").Value = -1324.0814
